$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "GFG"
$ws.Range("B13").Value = "Check If Circular Linked List"

# Match the style used by the other recently-added rows (wrap text), same as B10:B12
$ws.Range("B13").WrapText = $true

# Restore the view state captured in the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("H9").Select() | Out-Null
